$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "groups"/"counts" values in column A (rows 2-6) to their new figures
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 0

# Column B only changes on rows 3 and 6
$ws.Range("B3").Value = 4
$ws.Range("B6").Value = 2

# Rows 7 and 8 no longer exist in the updated table - remove them so the
# sheet (and its dimension) shrinks from A1:B8 down to A1:B6
$ws.Range("A7:B8").EntireRow.Delete()
